# Apply value updates to Sheet1 as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("I4").Value  = 1.73
$ws.Range("L4").Value  = 2.38
$ws.Range("M4").Value  = 1.07
$ws.Range("N4").Value  = 9
$ws.Range("U4").Value  = 2
$ws.Range("V4").Value  = 1.73
$ws.Range("W4").Value  = 11
$ws.Range("AE4").Value = 19
$ws.Range("AH4").Value = 7.5
$ws.Range("AQ4").Value = 101
$ws.Range("AS4").Value = 301
$ws.Range("AW4").Value = 3.6
$ws.Range("AX4").Value = 9
$ws.Range("AZ4").Value = 29

# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9

# Row 12
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3

# Row 17
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 15

# Row 20
$ws.Range("N20").Value = 8
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 1.65

# Row 26
$ws.Range("G26").Value  = 2.92
$ws.Range("I26").Value  = 2.65
$ws.Range("J26").Value  = 3.6
$ws.Range("Q26").Value  = 2.5
$ws.Range("R26").Value  = 1.4
$ws.Range("W26").Value  = 6.6
$ws.Range("X26").Value  = 13.5
$ws.Range("Z26").Value  = 37
$ws.Range("AG26").Value = 6.3
$ws.Range("AH26").Value = 11.75
$ws.Range("AJ26").Value = 32
$ws.Range("AK26").Value = 28
$ws.Range("AL26").Value = 45
$ws.Range("AN26").Value = 4.55
$ws.Range("AP26").Value = 29
$ws.Range("AW26").Value = 4.3

$wb.Save()
